$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.514.00"
$ws.Range("E2").Value = "  +0.16%  "

$ws.Range("D3").Value = "2.724.15"
$ws.Range("E3").Value = "  +2.91%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "610.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.71%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.553"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.33%  "

$ws.Range("D9").Value = "2.723.71"
$ws.Range("E9").Value = "  +2.92%  "

$ws.Range("E10").Value = "  +1.97%  "

$ws.Range("E11").Value = "  +4.51%  "

$ws.Range("E12").Value = "  +0.51%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.32"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.53%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.35%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000189"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.18%  "

$ws.Range("D17").Value = "68.498.06"
$ws.Range("E17").Value = "  +0.37%  "

$ws.Range("D18").Value = "2.683.91"
$ws.Range("E18").Value = "  +2.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "372.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.68%  "

$ws.Range("E21").Value = "  +3.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.66%  "

$ws.Range("E23").Value = "  +4.72%  "

$ws.Range("E24").Value = "  +1.31%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.57%  "

$ws.Range("E26").Value = "  +0.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.03%  "

$ws.Range("E28").Value = "  +3.08%  "

$ws.Range("E29").Value = "  +1.23%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "589.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.984"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.84%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.32%  "

$ws.Range("E33").Value = "  +3.07%  "

$ws.Range("E34").Value = "  +6.14%  "

$ws.Range("E35").Value = "  +2.41%  "

$ws.Range("E36").Value = "  -3.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "162.83"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.48%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.92"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.43%  "

$ws.Range("E40").Value = "  +2.61%  "

$ws.Range("E41").Value = "  +1.62%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.99"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.00%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.77%  "

$ws.Range("E45").Value = "  +0.00%  "

$ws.Range("E46").Value = "  -3.31%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.61%  "

$ws.Range("E48").Value = "  +4.40%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "155.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.71%  "

$ws.Range("E50").Value = "  +3.29%  "

$ws.Range("E51").Value = "  +5.29%  "
